$d = $word.ActiveDocument

$wordOpenXmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-PkgXml($bodyFragment) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document ' + $wordOpenXmlNs + '><w:body>' + $bodyFragment + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# -----------------------------------------------------------------
# Change 1: Ingredients list - "Crab-apples" line (paragraph 2).
# Replace <w:tab/> + "Crab-apples" + "\n" runs with a 16-space run
# (no formatting) followed by a single "Crab-apples\n" run.
# -----------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$p2Full = $p2.Range
$p2Inner = $d.Range($p2Full.Start, $p2Full.End - 1)
$p2Xml = New-PkgXml('<w:p><w:r><w:t xml:space="preserve">                </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>Crab-apples\n</w:t></w:r></w:p>')
$p2Inner.InsertXML($p2Xml)

# -----------------------------------------------------------------
# Change 2: Ingredients list - "lemon juice or vinegar" line (paragraph 3).
# Replace the leading <w:tab/> run with a 16-space run (no formatting);
# the rest of the paragraph is untouched.
# -----------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$p3Start = $p3.Range.Start
$p3ZeroRange = $d.Range($p3Start, $p3Start)
$p3SpacesXml = New-PkgXml('<w:p><w:r><w:t xml:space="preserve">                </w:t></w:r></w:p>')
$p3ZeroRange.InsertXML($p3SpacesXml)
# The original <w:tab/> run now sits right after the 16 inserted spaces.
$p3TabPos = $p3Start + 16
$p3TabRange = $d.Range($p3TabPos, $p3TabPos + 1)
$p3TabRange.Delete()

# -----------------------------------------------------------------
# Change 3: "Instructions" heading (paragraph 5).
# Demote it from the Heading2 style to a plain paragraph, keeping the
# text "Instructions\n".
# -----------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5)
$p5Full = $p5.Range
$p5Xml = New-PkgXml('<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000"/></w:rPr><w:t>Instructions\n</w:t></w:r></w:p>')
$p5Full.InsertXML($p5Xml)

# -----------------------------------------------------------------
# Change 4: "Slice the Crab-apples evenly into 1/8-inch-thick pieces."
# (paragraph 12). Merge the trailing " " run and the
# "evenly into 1/8-inch-thick pieces.\n" run into a single run.
# -----------------------------------------------------------------
$p12 = $d.Paragraphs.Item(12)
$p12Full = $p12.Range
$p12Xml = New-PkgXml('<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="720"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">Slice the </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>Crab-apples</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve"> evenly into 1/8-inch-thick pieces.\n</w:t></w:r></w:p>')
$p12Full.InsertXML($p12Xml)

Write-Output "done"
